$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K5").Value = 21.28240740740739
$ws.Range("R5").Value = 2.017497406510892
$ws.Range("S5").Value = 2.228623569098047

$ws.Range("K6").Value = 21.28240740740739
$ws.Range("R6").Value = 1.368667816644515
$ws.Range("S6").Value = 1.441442038370027

$ws.Range("K7").Value = 21.28240740740739

$ws.Range("K20").Value = 3.38888888888889
$ws.Range("R20").Value = 1.233817681248088
$ws.Range("S20").Value = 1.291146001942376

$ws.Range("K21").Value = 3.38888888888889
$ws.Range("R21").Value = 1.694051767048283
$ws.Range("S21").Value = 1.836167304537999

$ws.Range("K22").Value = 3.38888888888889

$ws.Range("K38").Value = -1.226851851851833
$ws.Range("R38").Value = 1.203236793039155
$ws.Range("S38").Value = 1.257328254301852

$ws.Range("K39").Value = -1.226851851851833
$ws.Range("R39").Value = 1.626775542720574
$ws.Range("S39").Value = 1.756382654173023

$ws.Range("K40").Value = -1.226851851851833

$ws.Range("K51").Value = 21.79166666666666
$ws.Range("R51").Value = 1.372938473321419
$ws.Range("S51").Value = 1.446233342398694

$ws.Range("K52").Value = 21.79166666666666

$ws.Range("K53").Value = 21.79166666666666
$ws.Range("R53").Value = 2.028520339740724
$ws.Range("S53").Value = 2.242263395092639

$ws.Range("K54").Value = 21.19907407407406
$ws.Range("R54").Value = 1.367971510132557
$ws.Range("S54").Value = 1.440661027663225

$ws.Range("K55").Value = 21.19907407407406
$ws.Range("R55").Value = 2.015705049109126
$ws.Range("S55").Value = 2.22640738080769

$ws.Range("K56").Value = 21.19907407407406

$ws.Range("K57").Value = 21.19907407407406
$ws.Range("R57").Value = 1.367971510132557
$ws.Range("S57").Value = 1.440661027663225

$ws.Range("K58").Value = 21.19907407407406
$ws.Range("R58").Value = 2.015705049109126
$ws.Range("S58").Value = 2.22640738080769

$ws.Range("K59").Value = 21.19907407407406

$ws.Range("K72").Value = -1.226851851851833

$ws.Range("K73").Value = -1.226851851851833
$ws.Range("R73").Value = 1.626775542720574
$ws.Range("S73").Value = 1.756382654173023

$ws.Range("K74").Value = -1.226851851851833
$ws.Range("R74").Value = 1.203236793039155
$ws.Range("S74").Value = 1.257328254301852

$ws.Range("K84").Value = 13.17361111111111
$ws.Range("R84").Value = 1.304077921028169
$ws.Range("S84").Value = 1.369214264257821

$ws.Range("K85").Value = 13.17361111111111

$ws.Range("K86").Value = 13.17361111111111

$ws.Range("K87").Value = 13.17361111111111

$ws.Range("K88").Value = 13.17361111111111

$ws.Range("K89").Value = 13.17361111111111
$ws.Range("R89").Value = 1.856836936506854
$ws.Range("S89").Value = 2.031823338122968

$ws.Range("K90").Value = 13.17361111111111
